$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2023 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 300000000
$ws.Cells.Item($row, 7).Value = "Espárragos"
$ws.Cells.Item($row, 8).Value = "Verde"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 1200
$ws.Cells.Item($row, 11).Value = 2000
$ws.Cells.Item($row, 12).Value = 2000
$ws.Cells.Item($row, 13).Value = 2000
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Linares"
$ws.Cells.Item($row, 16).Value = 2000
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
